# edit.ps1 - applies the content changes described by the commit diff
# (style-id renumbering / namespace shuffling in the diff are artifacts of
# which application last saved the package and are not reproduced here;
# only the genuine textual edits are applied.)

$d = $word.ActiveDocument

# 1) "...saving 80% of the work for the following years" -> "...years" -> "...year"
$d.Content.Find.Execute(
    "for the following years",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "for the following year",
    2
) | Out-Null

# 2) Rewrite part of the Radar bullet:
#    "...and revealed potential implications and ramifications by documenting
#     with completeness, accuracy, and consistency perspectives "
# -> "...and documented from completeness, accuracy, and consistency
#     perspectives to reveal potential implications and ramifications"
$d.Content.Find.Execute(
    "and revealed potential implications and ramifications by documenting with completeness, accuracy",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and documented from completeness, accuracy",
    2
) | Out-Null

$d.Content.Find.Execute(
    "consistency perspectives ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "consistency perspectives to reveal potential implications and ramifications",
    2
) | Out-Null
